$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$ws.Range("B11").Value = "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."
$ws.Range("C11").Value = "This discipline aims to present the fundamentals of Thermodynamics and Statistical Physics. To present the different formalisms of Statistical Physics. Application of formalisms to some simple models. Applications."

$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

$ws.Range("B14").Value = "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."
$ws.Range("C14").Value = "Concepts of thermodynamics. Formalisms of Statistical Physics. Classic ideal gas. Quantum gases. Applications."

$ws.Range("B15").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C15").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$ws.Range("B16").Value = "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."
$ws.Range("C16").Value = "Macroscopic and microscopic systems. Postulates of thermodynamics. Fundamental equation. Equations of state. Thermodynamic equilibrium. Thermodynamic derivatives. Thermodynamic potentials. Maxwell relations. Born diagram. Reduction of thermodynamic derivatives.Microcanonical formalism. Boltzmann equation. Einstein model of a crystalline solid. Canonical formalism. Example: ideal classical gas. Maxwell-Boltzmann distribution. Grand canonical formalism. Quantum gases. Fermions and bosons. Bose-Einstein distribution. Fermi-Dirac distribution. Examples: electron gas and photon gas. Planck distribution.Applications: superconductivity, electron gas in semiconductor, superfluidity of the liquid helium."
